$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") values must be written as text (not numbers) so that
# formatting such as trailing zeros, thousands separators written with dots,
# etc. is preserved exactly like in the source data. Excel auto-detects
# numeric-looking strings as numbers unless the cell is pre-formatted as
# text, so set NumberFormat="@" on each Price cell individually first
# (union/multi-area ranges are not reliably honored), write the value, then
# restore the default "Normal" style so no stray formatting is left behind.
$priceCells = @("D2","D3","D4","D5","D6","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D27","D28","D29","D30","D31","D33","D35","D36","D37","D38","D40","D41","D42","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '58.780.25'
$ws.Range('E2').Value = '  -4.62%  '
$ws.Range('D3').Value = '2.469.55'
$ws.Range('E3').Value = '  -4.10%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '532.95'
$ws.Range('E5').Value = '  -3.48%  '
$ws.Range('D6').Value = '144.04'
$ws.Range('E6').Value = '  -6.39%  '
$ws.Range('E7').Value = '  -0.26%  '
$ws.Range('D8').Value = '0.569'
$ws.Range('E8').Value = '  -4.41%  '
$ws.Range('D9').Value = '2.492.04'
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('D10').Value = '0.0995'
$ws.Range('E10').Value = '  -4.28%  '
$ws.Range('E11').Value = '  -2.16%  '
$ws.Range('D12').Value = '5.61'
$ws.Range('E12').Value = '  +3.32%  '
$ws.Range('D13').Value = '0.351'
$ws.Range('E13').Value = '  -3.32%  '
$ws.Range('D14').Value = '2.903.38'
$ws.Range('E14').Value = '  -4.21%  '
$ws.Range('D15').Value = '23.81'
$ws.Range('E15').Value = '  -5.98%  '
$ws.Range('D16').Value = '58.623.05'
$ws.Range('E16').Value = '  -4.70%  '
$ws.Range('D17').Value = '0.0000138'
$ws.Range('E17').Value = '  -3.81%  '
$ws.Range('D18').Value = '2.475.41'
$ws.Range('E18').Value = '  -4.02%  '
$ws.Range('D19').Value = '11.30'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('D20').Value = '4.30'
$ws.Range('E20').Value = '  -4.95%  '
$ws.Range('D21').Value = '322.26'
$ws.Range('E21').Value = '  -4.24%  '
$ws.Range('D22').Value = '0.997'
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('D23').Value = '5.73'
$ws.Range('E23').Value = '  -4.92%  '
$ws.Range('D24').Value = '60.66'
$ws.Range('E24').Value = '  -3.58%  '
$ws.Range('E25').Value = '  -11.37%  '
$ws.Range('E26').Value = '  -0.25%  '
$ws.Range('D27').Value = '0.161'
$ws.Range('E27').Value = '  -4.36%  '
$ws.Range('D28').Value = '2.579.67'
$ws.Range('E28').Value = '  -4.48%  '
$ws.Range('D29').Value = '7.71'
$ws.Range('E29').Value = '  -3.87%  '
$ws.Range('D30').Value = '6.91'
$ws.Range('E30').Value = '  -3.99%  '
$ws.Range('D31').Value = '0.0₃0772'
$ws.Range('E31').Value = '  -7.32%  '
$ws.Range('E32').Value = '  -6.82%  '
$ws.Range('D33').Value = '1.24'
$ws.Range('E33').Value = '  -4.68%  '
$ws.Range('E34').Value = '  -0.25%  '
$ws.Range('D35').Value = '157.84'
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').Value = '1.40'
$ws.Range('E36').Value = '  -0.73%  '
$ws.Range('D37').Value = '18.51'
$ws.Range('E37').Value = '  -3.36%  '
$ws.Range('D38').Value = '4.38'
$ws.Range('E38').Value = '  -5.44%  '
$ws.Range('E39').Value = '  -10.13%  '
$ws.Range('D40').Value = '5.74'
$ws.Range('E40').Value = '  -4.14%  '
$ws.Range('D41').Value = '304.79'
$ws.Range('E41').Value = '  -8.46%  '
$ws.Range('D42').Value = '36.52'
$ws.Range('E42').Value = '  -2.38%  '
$ws.Range('D43').Value = '3.71'
$ws.Range('E43').Value = '  -5.72%  '
$ws.Range('D44').Value = '0.807'
$ws.Range('E44').Value = '  -10.32%  '
$ws.Range('E45').Value = '  -0.27%  '
$ws.Range('D47').Value = '0.592'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('D48').Value = '124.13'
$ws.Range('E48').Value = '  +0.09%  '
$ws.Range('D49').Value = '0.0924'
$ws.Range('E49').Value = '  -4.33%  '
$ws.Range('D50').Value = '0.0519'
$ws.Range('E50').Value = '  -4.74%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '18.50'
$ws.Range('E51').Value = '  -5.23%  '

foreach ($addr in $priceCells) {
    $ws.Range($addr).Style = "Normal"
}

Write-Host "Updated cryptos list"